$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.4557103333333333
$ws.Cells.Item(2, 8).Value = 1.367131
$ws.Cells.Item(2, 9).Value = 0.1996401272959883
$ws.Cells.Item(2, 10).Value = 0.1996401272959883
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.4884036666666667
$ws.Cells.Item(2, 14).Value = 1.465211
$ws.Cells.Item(2, 15).Value = 0.03210371182687088
$ws.Cells.Item(2, 16).Value = 0.03210371182687088
$ws.Cells.Item(2, 17).Value = 0.2225705977378889
$ws.Cells.Item(2, 18).Value = 2.003135379641
$ws.Cells.Item(2, 19).Value = 0.006409189115790228
$ws.Cells.Item(2, 20).Value = 0.006409189115790226
$ws.Cells.Item(3, 7).Value = 0.4557103333333333
$ws.Cells.Item(3, 8).Value = 1.367131
$ws.Cells.Item(3, 9).Value = 0.1996401272959883
$ws.Cells.Item(3, 10).Value = 0.1996401272959883
$ws.Cells.Item(3, 15).Value = 0.9336784186214153
$ws.Cells.Item(3, 16).Value = 0.9336784186214153
$ws.Cells.Item(3, 17).Value = 6.473063452855889
$ws.Cells.Item(3, 18).Value = 58.257571075703
$ws.Cells.Item(3, 19).Value = 0.1863996783470964
$ws.Cells.Item(3, 20).Value = 0.1863996783470964
$ws.Cells.Item(4, 7).Value = 0.4557103333333333
$ws.Cells.Item(4, 8).Value = 1.367131
$ws.Cells.Item(4, 9).Value = 0.1996401272959883
$ws.Cells.Item(4, 10).Value = 0.1996401272959883
$ws.Cells.Item(4, 13).Value = 0.520567
$ws.Cells.Item(4, 14).Value = 1.561701
$ws.Cells.Item(4, 15).Value = 0.03421786955171377
$ws.Cells.Item(4, 16).Value = 0.03421786955171376
$ws.Cells.Item(4, 17).Value = 0.2372277610923333
$ws.Cells.Item(4, 18).Value = 2.135049849831
$ws.Cells.Item(4, 19).Value = 0.006831259833101659
$ws.Cells.Item(4, 20).Value = 0.006831259833101656
$ws.Cells.Item(5, 9).Value = 0.2962807848215612
$ws.Cells.Item(5, 10).Value = 0.2962807848215612
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.4884036666666667
$ws.Cells.Item(5, 14).Value = 1.465211
$ws.Cells.Item(5, 15).Value = 0.03210371182687088
$ws.Cells.Item(5, 16).Value = 0.03210371182687088
$ws.Cells.Item(5, 17).Value = 0.330311306996
$ws.Cells.Item(5, 18).Value = 2.972801762964
$ws.Cells.Item(5, 19).Value = 0.009511712935750543
$ws.Cells.Item(5, 20).Value = 0.009511712935750538
$ws.Cells.Item(6, 9).Value = 0.2962807848215612
$ws.Cells.Item(6, 10).Value = 0.2962807848215612
$ws.Cells.Item(6, 15).Value = 0.9336784186214153
$ws.Cells.Item(6, 16).Value = 0.9336784186214153
$ws.Cells.Item(6, 19).Value = 0.2766309746401071
$ws.Cells.Item(6, 20).Value = 0.2766309746401071
$ws.Cells.Item(7, 9).Value = 0.2962807848215612
$ws.Cells.Item(7, 10).Value = 0.2962807848215612
$ws.Cells.Item(7, 13).Value = 0.520567
$ws.Cells.Item(7, 14).Value = 1.561701
$ws.Cells.Item(7, 15).Value = 0.03421786955171377
$ws.Cells.Item(7, 16).Value = 0.03421786955171376
$ws.Cells.Item(7, 17).Value = 0.352063626636
$ws.Cells.Item(7, 18).Value = 3.168572639724
$ws.Cells.Item(7, 19).Value = 0.01013809724570356
$ws.Cells.Item(7, 20).Value = 0.01013809724570356
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.3513206666666667
$ws.Cells.Item(8, 8).Value = 1.053962
$ws.Cells.Item(8, 9).Value = 0.1539085192605057
$ws.Cells.Item(8, 10).Value = 0.1539085192605057
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.4884036666666667
$ws.Cells.Item(8, 14).Value = 1.465211
$ws.Cells.Item(8, 15).Value = 0.03210371182687088
$ws.Cells.Item(8, 16).Value = 0.03210371182687088
$ws.Cells.Item(8, 17).Value = 0.1715863017757778
$ws.Cells.Item(8, 18).Value = 1.544276715982
$ws.Cells.Item(8, 19).Value = 0.004941034750039683
$ws.Cells.Item(8, 20).Value = 0.004941034750039681
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.3513206666666667
$ws.Cells.Item(9, 8).Value = 1.053962
$ws.Cells.Item(9, 9).Value = 0.1539085192605057
$ws.Cells.Item(9, 10).Value = 0.1539085192605057
$ws.Cells.Item(9, 15).Value = 0.9336784186214153
$ws.Cells.Item(9, 16).Value = 0.9336784186214153
$ws.Cells.Item(9, 17).Value = 4.990277378611778
$ws.Cells.Item(9, 18).Value = 44.912496407506
$ws.Cells.Item(9, 19).Value = 0.1437010628755126
$ws.Cells.Item(9, 20).Value = 0.1437010628755126
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.3513206666666667
$ws.Cells.Item(10, 8).Value = 1.053962
$ws.Cells.Item(10, 9).Value = 0.1539085192605057
$ws.Cells.Item(10, 10).Value = 0.1539085192605057
$ws.Cells.Item(10, 13).Value = 0.520567
$ws.Cells.Item(10, 14).Value = 1.561701
$ws.Cells.Item(10, 15).Value = 0.03421786955171377
$ws.Cells.Item(10, 16).Value = 0.03421786955171376
$ws.Cells.Item(10, 17).Value = 0.1828859454846667
$ws.Cells.Item(10, 18).Value = 1.645973509362
$ws.Cells.Item(10, 19).Value = 0.005266421634953411
$ws.Cells.Item(10, 20).Value = 0.005266421634953409
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.3449053333333333
$ws.Cells.Item(11, 8).Value = 1.034716
$ws.Cells.Item(11, 9).Value = 0.1510980542136751
$ws.Cells.Item(11, 10).Value = 0.1510980542136751
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.4884036666666667
$ws.Cells.Item(11, 14).Value = 1.465211
$ws.Cells.Item(11, 15).Value = 0.03210371182687088
$ws.Cells.Item(11, 16).Value = 0.03210371182687088
$ws.Cells.Item(11, 17).Value = 0.1684530294528889
$ws.Cells.Item(11, 18).Value = 1.516077265076
$ws.Cells.Item(11, 19).Value = 0.004850808390076739
$ws.Cells.Item(11, 20).Value = 0.004850808390076737
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.3449053333333333
$ws.Cells.Item(12, 8).Value = 1.034716
$ws.Cells.Item(12, 9).Value = 0.1510980542136751
$ws.Cells.Item(12, 10).Value = 0.1510980542136751
$ws.Cells.Item(12, 15).Value = 0.9336784186214153
$ws.Cells.Item(12, 16).Value = 0.9336784186214153
$ws.Cells.Item(12, 17).Value = 4.899151817700889
$ws.Cells.Item(12, 18).Value = 44.092366359308
$ws.Cells.Item(12, 19).Value = 0.141076992314997
$ws.Cells.Item(12, 20).Value = 0.141076992314997
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.3449053333333333
$ws.Cells.Item(13, 8).Value = 1.034716
$ws.Cells.Item(13, 9).Value = 0.1510980542136751
$ws.Cells.Item(13, 10).Value = 0.1510980542136751
$ws.Cells.Item(13, 13).Value = 0.520567
$ws.Cells.Item(13, 14).Value = 1.561701
$ws.Cells.Item(13, 15).Value = 0.03421786955171377
$ws.Cells.Item(13, 16).Value = 0.03421786955171376
$ws.Cells.Item(13, 17).Value = 0.1795463346573333
$ws.Cells.Item(13, 18).Value = 1.615917011916
$ws.Cells.Item(13, 19).Value = 0.005170253508601309
$ws.Cells.Item(13, 20).Value = 0.005170253508601307
$ws.Cells.Item(14, 7).Value = 0.4544146666666666
$ws.Cells.Item(14, 8).Value = 1.363244
$ws.Cells.Item(14, 9).Value = 0.1990725144082698
$ws.Cells.Item(14, 10).Value = 0.1990725144082698
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.4884036666666667
$ws.Cells.Item(14, 14).Value = 1.465211
$ws.Cells.Item(14, 15).Value = 0.03210371182687088
$ws.Cells.Item(14, 16).Value = 0.03210371182687088
$ws.Cells.Item(14, 17).Value = 0.2219377893871111
$ws.Cells.Item(14, 18).Value = 1.997440104484
$ws.Cells.Item(14, 19).Value = 0.006390966635213696
$ws.Cells.Item(14, 20).Value = 0.006390966635213693
$ws.Cells.Item(15, 7).Value = 0.4544146666666666
$ws.Cells.Item(15, 8).Value = 1.363244
$ws.Cells.Item(15, 9).Value = 0.1990725144082698
$ws.Cells.Item(15, 10).Value = 0.1990725144082698
$ws.Cells.Item(15, 15).Value = 0.9336784186214153
$ws.Cells.Item(15, 16).Value = 0.9336784186214153
$ws.Cells.Item(15, 17).Value = 6.454659366019111
$ws.Cells.Item(15, 18).Value = 58.091934294172
$ws.Cells.Item(15, 19).Value = 0.1858697104437022
$ws.Cells.Item(15, 20).Value = 0.1858697104437022
$ws.Cells.Item(16, 7).Value = 0.4544146666666666
$ws.Cells.Item(16, 8).Value = 1.363244
$ws.Cells.Item(16, 9).Value = 0.1990725144082698
$ws.Cells.Item(16, 10).Value = 0.1990725144082698
$ws.Cells.Item(16, 13).Value = 0.520567
$ws.Cells.Item(16, 14).Value = 1.561701
$ws.Cells.Item(16, 15).Value = 0.03421786955171377
$ws.Cells.Item(16, 16).Value = 0.03421786955171376
$ws.Cells.Item(16, 17).Value = 0.2365532797826667
$ws.Cells.Item(16, 18).Value = 2.128979518044
$ws.Cells.Item(16, 19).Value = 0.006811837329353835
$ws.Cells.Item(16, 20).Value = 0.006811837329353833
